$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Пыакео сукчалын" + <space> + "Группа:"  ->  "Пыакео сукчалын" + <br/>
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Пыакео сукчалын", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRng = $d.Range($rng.End, $rng.End + 1)
$spaceRng.InsertBreak(6)

# ---------------------------------------------------------------------------
# 2) Rename bookmark "создание-отчета" -> "компиляция-отчета-с-помощью-make"
#    and update its heading text.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("создание-отчета")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("компиляция-отчета-с-помощью-make", $bmRange)

$d.Content.Find.Execute("3. Создание отчета", $true, $false, $false, $false, $false, $true, 1, $false, "3. Компиляция отчета с помощью make", 2)

# ---------------------------------------------------------------------------
# 3) Replace the "Создан отчет..." paragraph with the new text and append
#    four new paragraphs (two Heading4 + two FirstParagraph) using the
#    paragraph-mark find/replace code ^p so each chunk becomes its own
#    <w:p> inheriting the FirstParagraph style from the source paragraph.
# ---------------------------------------------------------------------------
$oldReport = "Создан отчет с использованием различных элементов Markdown."
$newReport = "Проведена компиляция отчета с использованием утилиты make.^p4. Проверка сгенерированных файлов^pПосле компиляции созданы файлы report.pdf и report.docx.^p5. Загрузка на GitHub^pОтчет успешно загружен на GitHub."
$d.Content.Find.Execute($oldReport, $true, $false, $false, $false, $false, $true, 1, $false, $newReport, 2)

# Re-find the two new heading paragraphs and promote them to Heading4,
# bookmarking each exactly like the existing headings in the document.
$f = $d.Content.Find
$f.Execute("4. Проверка сгенерированных файлов", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $f.Duplicate.Paragraphs(1)
$headingPara.Style = "Heading4"
$hRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$d.Bookmarks.Add("проверка-сгенерированных-файлов", $hRange)

$f2 = $d.Content.Find
$f2.Execute("5. Загрузка на GitHub", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara2 = $f2.Duplicate.Paragraphs(1)
$headingPara2.Style = "Heading4"
$hRange2 = $d.Range($headingPara2.Range.Start, $headingPara2.Range.End - 1)
$d.Bookmarks.Add("загрузка-на-github", $hRange2)

# ---------------------------------------------------------------------------
# 4) Insert a new BodyText paragraph (bold label run + space run + url run)
#    right after "Оформлен отчет по лабораторной работе №2 в формате Markdown."
# ---------------------------------------------------------------------------
$f3 = $d.Content.Find
$f3.Execute("Оформлен отчет по лабораторной работе №2 в формате Markdown.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$srcPara = $f3.Duplicate.Paragraphs(1)
$srcPara.Range.InsertParagraphAfter()
$linkPara = $srcPara.Next()
$linkPara.Style = "BodyText"

$label = "Ссылка на репозиторий GitHub:"
$url = "https://github.com/soukchalern3-blip/study_2025-2026_arch-pc"
$linkPara.Range.Text = $label + " " + $url

# Drop a throwaway bookmark right at the space/url boundary so the engine
# keeps that boundary as a separate run instead of coalescing it back into
# the run that follows once formatting is applied; then bold just the
# label range (must be the very last text/format touch in this paragraph
# so the bold attribute does not leak into the later runs).
$splitPos = $linkPara.Range.Start + ($label + " ").Length
$tmpBmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TEMP_SPLIT_MARK", $tmpBmRange)

$labelEnd = $linkPara.Range.Start + $label.Length
$labelRng = $d.Range($linkPara.Range.Start, $labelEnd)
$labelRng.Bold = 1

$d.Bookmarks("TEMP_SPLIT_MARK").Delete()

# ---------------------------------------------------------------------------
# 5) Merge the split "...отче" + " " + "тов..." runs in the Выводы section.
# ---------------------------------------------------------------------------
$oldConclusion = "В ходе работы были освоены основные принципы работы с языком разметки Markdown. Приобретены навыки оформления отче тов с использованием данного инструмента."
$newConclusion = "В ходе работы были освоены основные принципы работы с языком разметки Markdown. Приобретены навыки оформления отчетов с использованием данного инструмента."
$d.Content.Find.Execute($oldConclusion, $true, $false, $false, $false, $false, $true, 1, $false, $newConclusion, 2)

Write-Host "edit complete"
